$d = $word.ActiveDocument
$last = $d.Paragraphs.Last
$rng = $last.Range
$rng.Collapse(1)

$body = @'
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>chances are there will be if</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>احتمال وجود خواهد داشت اگر</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">virtuous circle: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">mend: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>رفو کردن، تعمیر کردن لباس</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>Intimidate:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> مرعو</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="eastAsia"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ب</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve"> کردن ، ترساندن </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Profile: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>نیم رخ، شرح حال، شکل دادن</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Horrendous: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>دهشتناک</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Subcontinent: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>شبه قاره</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Serialize: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>مرتب کردن</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Crude plan: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>طرح خام، طرح ناپخته</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Stark: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>کامل، رک و سرراست</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Perpetual: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>رایج، مرسوم</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Peculiar: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ابدی، همیشگی</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Anonymity: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>ناشناس</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Vie: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>رقابت کردن، هم چشمی کردن</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">Mending </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:ind w:left="900"/>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>3</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t xml:space="preserve">clothes: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:cs="B Compset" w:hint="cs"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:rtl/>
          <w:lang w:bidi="fa-IR"/>
        </w:rPr>
        <w:t>وصله کردن، رفو کردن لباس</w:t>
      </w:r>
    </w:p>
'@

$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($frag)

Write-Host "Paragraph count after: $($d.Paragraphs.Count)"
